# Daily Report template update:
#  - Row 2 / Row 3 "Activity" and "Result / Actions" text replaced with new
#    content about learning high-level Python features/functions.
#  - Row 2 height reverts to the default (shorter text), row 3 grows to the
#    taller (27pt) height that row 2 used to have (longer text).
#  - Selection moves from A1:C3 to A4:C4 (active cell A4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the two data rows with the new activity / result text ---
$ws.Range("B2").Value = "learned high-level features"
$ws.Range("C2").Value = "Including iteration(Iterable, Iterator) and generator, list generator"

$ws.Range("B3").Value = "learned high-level function "
$ws.Range("C3").Value = "Including return function(return value is a function), lambda function, decorator, partial function, high-level built-in functions like map/reduce and filter, sorted"

# --- Row heights follow the new (shorter/longer) wrapped text ---
# Row 2's text got shorter, so it shrinks back to the default row height.
$ws.Rows.Item(2).AutoFit()
# Row 3's text got longer (now wraps to two lines), so it grows.
$ws.Rows.Item(3).RowHeight = 27

# --- Selection moves to the empty merged row beneath the table ---
$ws.Range("A4:C4").Select()
